$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" (Total) summary sheet: insert a new 2022-Q4 row at the top of the
#    data, pushing the existing 2022-Q3 / 2022-Q2 rows down by one.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Push old row 3 (2022-Q2) down into row 4, carrying the row-2/3 formatting
# (bold/border style) along with it so the new A4 matches A2/A3.
$total.Range("A3").Copy($total.Range("A4"))
$total.Range("A4").Value = 2
$total.Range("B4").Value = $total.Range("B3").Value2
$total.Range("C4").Value = $total.Range("C3").Value2
$total.Range("D4").Value = $total.Range("D3").Value2

# Push old row 2 (2022-Q3) down into row 3.
$total.Range("A3").Value = 1
$total.Range("B3").Value = $total.Range("B2").Value2
$total.Range("C3").Value = $total.Range("C2").Value2
$total.Range("D3").Value = $total.Range("D2").Value2

# Write the brand-new 2022-Q4 summary row into row 2.
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.04

# ---------------------------------------------------------------------------
# 2) Add a brand-new "2022-Q4" fund-holdings sheet (same shape/style as the
#    existing "2022-Q3" sheet), placed right after "总计" and before the old
#    "2022-Q3" sheet. It carries the refreshed fund figures.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

$q4.Range("D2").Value = "'0.90"
$q4.Range("E2").Value = "'94.26"
$q4.Range("F2").Value = "'3.43"
$q4.Range("G2").Value = "'0.0309"
$q4.Range("H2").Value = 7

$q4.Range("D3").Value = "'0.20"
$q4.Range("E3").Value = "'94.26"
$q4.Range("F3").Value = "'3.43"
$q4.Range("G3").Value = "'0.0069"
$q4.Range("H3").Value = 7
